# The "reviews_count" column (column E) is removed from the sheet.
# Deleting the entire column shifts every column to its right
# (reviews_average, latitude, longitude, is_permanently_closed,
# gmaps_link, latest_review_date) one position to the left, which is
# exactly the transformation described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(5).Delete()
